$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (A4) so the new date cell (A5)
# gets the same date number format/style instead of creating a new style.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row of data (row 5)
$ws.Range("A5").Value = 44314
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "Сделал guards-ы, страницу просмотра заявки, получение данных с стороннего REST api."

# Update the active selection to the newly added cell, like the diff shows
$ws.Range("C5").Select() | Out-Null

# Recalculate so the SUM formula in F2 reflects the newly added hours
$excel.Calculate() | Out-Null
